$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells whose new values would otherwise
# be auto-converted to numbers by Excel (losing exact decimal formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values from the latest crypto data refresh
$ws.Range("D2").Value = "61.242.62"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "2.969.15"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "527.65"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "130.01"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "2.965.24"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "6.10"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "33.07"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").Value = "3.455.40"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "61.347.61"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.971.36"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "455.31"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "13.04"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "6.78"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").Value = "77.05"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "11.71"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").Value = "7.57"
$ws.Range("E28").Value = "  -6.99%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "25.29"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "55.29"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.31"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").Value = "2.22"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").Value = "5.74"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").Value = "446.84"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("D38").Value = "3.121.42"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "0.0380"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "0.115"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").Value = "7.91"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  -5.99%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "24.65"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").Value = "119.86"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "1.91"
$ws.Range("E49").Value = "  -4.34%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0499"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "1.22"
$ws.Range("E51").Value = "  +4.80%  "
